# Crypto price/volume refresh - automated data pull (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.873.16'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '3.502.67'
$ws.Range("E3").Value = '  -1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.40'
$ws.Range("E5").Value = '  -1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.53'
$ws.Range("E6").Value = '  -1.31%  '

$ws.Range("D7").Value = '3.500.09'
$ws.Range("E7").Value = '  -1.20%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  +0.73%  '

$ws.Range("E10").Value = '  +3.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.56'
$ws.Range("E11").Value = '  +6.91%  '

$ws.Range("E13").Value = '  -1.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.02'
$ws.Range("E14").Value = '  +0.13%  '

$ws.Range("D15").Value = '4.095.60'
$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("D16").Value = '67.773.26'
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").Value = '3.507.17'
$ws.Range("E17").Value = '  -1.15%  '

$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("E19").Value = '  +1.25%  '

$ws.Range("E20").Value = '  +1.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.95'
$ws.Range("E21").Value = '  +2.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '445.44'
$ws.Range("E22").Value = '  -0.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.624'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.99'
$ws.Range("E24").Value = '  +2.57%  '

$ws.Range("D25").Value = '3.642.63'
$ws.Range("E25").Value = '  -1.16%  '

$ws.Range("E27").Value = '  -3.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.66'
$ws.Range("E28").Value = '  -4.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.97'
$ws.Range("E29").Value = '  -2.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.65'
$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("E31").Value = '  -1.65%  '

$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("E33").Value = '  +1.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.60'
$ws.Range("E35").Value = '  -0.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.17'
$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("E37").Value = '  -0.02%  '

$ws.Range("D38").Value = '3.498.32'
$ws.Range("E38").Value = '  -0.93%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.98'
$ws.Range("E39").Value = '  -0.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +6.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '176.30'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.895'
$ws.Range("E46").Value = '  +0.97%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.32'
$ws.Range("E47").Value = '  +5.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.73'
$ws.Range("E48").Value = '  +2.55%  '

$ws.Range("E49").Value = '  +0.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.53'
$ws.Range("E50").Value = '  -5.34%  '

$ws.Range("E51").Value = '  +0.19%  '

